$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Zanahoria at Vega Modelo de Temuco.
# It belongs chronologically among the existing records, so insert a fresh
# row at 214 (this pushes the former rows 214-225 down to 215-226) and fill
# it in with the new observation.
$ws.Rows(214).Insert()

$ws.Cells.Item(214, 1).Value = 10
$ws.Cells.Item(214, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(214, 3).Value = "La Araucanía"
$ws.Cells.Item(214, 4).Value = 44516
$ws.Cells.Item(214, 5).Value = 9
$ws.Cells.Item(214, 6).Value = 100114013
$ws.Cells.Item(214, 7).Value = "Zanahoria"
$ws.Cells.Item(214, 8).Value = "Sin especificar"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 55
$ws.Cells.Item(214, 11).Value = 10000
$ws.Cells.Item(214, 12).Value = 10000
$ws.Cells.Item(214, 13).Value = 10000
$ws.Cells.Item(214, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(214, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(214, 16).Value = 500
$ws.Cells.Item(214, 17).Value = 20
$ws.Cells.Item(214, 18).Value = "Hortaliza"
